# Update "想去人数" (want-to-go count) values across sheets to match latest
# generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 12
$ws1.Range("F12").Value = 227
$ws1.Range("F20").Value = 375
$ws1.Range("F23").Value = 19
$ws1.Range("F29").Value = 251
$ws1.Range("F31").Value = 569
$ws1.Range("F35").Value = 718
$ws1.Range("F36").Value = 95713

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F20").Value = 1721
$ws2.Range("F21").Value = 1721
$ws2.Range("F22").Value = 1110
$ws2.Range("F24").Value = 692
$ws2.Range("F26").Value = 14
$ws2.Range("F36").Value = 195

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 204
$ws3.Range("F7").Value = 4181
$ws3.Range("F10").Value = 371
$ws3.Range("F11").Value = 199
$ws3.Range("F12").Value = 221

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 204
$ws4.Range("F7").Value = 4181
$ws4.Range("F9").Value = 371
$ws4.Range("F10").Value = 371
$ws4.Range("F25").Value = 375
$ws4.Range("F27").Value = 1721
$ws4.Range("F28").Value = 1110
$ws4.Range("F32").Value = 251
$ws4.Range("F34").Value = 569
